$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 ("I0") and J1 ("IF")
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Match the header formatting used by the rest of row 1 (bold, centered,
# top-aligned, bordered) by copying H1's format onto the new header cells.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data values for columns I and J, rows 2-23
$values = @(
    @(2, 10, 10),
    @(3, 6, 6),
    @(4, 9, 9),
    @(5, 7, 7),
    @(6, 8, 8),
    @(7, 10, 10),
    @(8, 6, 6),
    @(9, 9, 9),
    @(10, 8, 8),
    @(11, 4, 4),
    @(12, 6, 6),
    @(13, 7, 7),
    @(14, 8, 8),
    @(15, 7, 7),
    @(16, 8, 8),
    @(17, 8, 8),
    @(18, 8, 9),
    @(19, 8, 8),
    @(20, 9, 9),
    @(21, 8, 8),
    @(22, 7, 7),
    @(23, 6, 6)
)

foreach ($row in $values) {
    $r = $row[0]
    $ws.Cells.Item($r, 9).Value = $row[1]
    $ws.Cells.Item($r, 10).Value = $row[2]
}
